$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force text format on all Price (D) cells being updated, so numeric-looking
# strings like "5.120" or "30.878.66" are preserved exactly as text.
foreach ($addr in @("D2","D3","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D36","D37","D38","D39","D40","D41","D42","D43","D45","D46","D47","D48","D49","D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '30.878.66'
$ws.Range("E2").Value = '  +2.58%  '
$ws.Range("D3").Value = '1.902.76'
$ws.Range("E3").Value = '  +1.07%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '246.36'
$ws.Range("E5").Value = '  +1.26%  '
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").Value = '0.5007'
$ws.Range("E7").Value = '  +0.96%  '
$ws.Range("D8").Value = '0.2995'
$ws.Range("E8").Value = '  +2.50%  '
$ws.Range("D9").Value = '0.06859'
$ws.Range("E9").Value = '  +3.80%  '
$ws.Range("D10").Value = '1.903.22'
$ws.Range("E10").Value = '  +1.15%  '
$ws.Range("D11").Value = '17.24'
$ws.Range("E11").Value = '  +2.55%  '
$ws.Range("D12").Value = '0.07345'
$ws.Range("E12").Value = '  +1.98%  '
$ws.Range("D13").Value = '91.82'
$ws.Range("E13").Value = '  +7.21%  '
$ws.Range("D14").Value = '5.120'
$ws.Range("E14").Value = '  +5.91%  '
$ws.Range("D15").Value = '0.6811'
$ws.Range("E15").Value = '  +2.82%  '
$ws.Range("D16").Value = '30.852.35'
$ws.Range("E16").Value = '  +2.55%  '
$ws.Range("D17").Value = '0.000008053'
$ws.Range("D18").Value = '13.35'
$ws.Range("E18").Value = '  +4.43%  '
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").Value = '2.154.84'
$ws.Range("E20").Value = '  +1.58%  '
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("D22").Value = '4.882'
$ws.Range("E22").Value = '  +2.47%  '
$ws.Range("D23").Value = '184.87'
$ws.Range("E23").Value = '  +36.76%  '
$ws.Range("D24").Value = '6.130'
$ws.Range("E24").Value = '  +9.45%  '
$ws.Range("D25").Value = '9.394'
$ws.Range("E25").Value = '  +2.92%  '
$ws.Range("D26").Value = '154.18'
$ws.Range("E26").Value = '  +1.56%  '
$ws.Range("D27").Value = '18.69'
$ws.Range("E27").Value = '  +11.50%  '
$ws.Range("D28").Value = '1.951'
$ws.Range("E28").Value = '  +2.17%  '
$ws.Range("D29").Value = '1.399'
$ws.Range("E29").Value = '  +1.19%  '
$ws.Range("D30").Value = '4.387'
$ws.Range("E30").Value = '  +5.50%  '
$ws.Range("D31").Value = '0.08986'
$ws.Range("E31").Value = '  +3.56%  '
$ws.Range("D32").Value = '4.085'
$ws.Range("E32").Value = '  +3.76%  '
$ws.Range("D33").Value = '0.05289'
$ws.Range("E33").Value = '  +5.97%  '
$ws.Range("D34").Value = '0.7455'
$ws.Range("E34").Value = '  +5.28%  '
$ws.Range("E35").Value = '  +3.46%  '
$ws.Range("D36").Value = '2.668'
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("D37").Value = '0.01931'
$ws.Range("E37").Value = '  +17.25%  '
$ws.Range("D38").Value = '2.729'
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("D39").Value = '2.193'
$ws.Range("E39").Value = '  +0.34%  '
$ws.Range("D40").Value = '0.9417'
$ws.Range("E40").Value = '  +0.77%  '
$ws.Range("D41").Value = '0.4411'
$ws.Range("E41").Value = '  +5.39%  '
$ws.Range("D42").Value = '106.22'
$ws.Range("E42").Value = '  +3.75%  '
$ws.Range("D43").Value = '5.855'
$ws.Range("E43").Value = '  -1.74%  '
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").Value = '7.785'
$ws.Range("E45").Value = '  +4.05%  '
$ws.Range("D46").Value = '0.1358'
$ws.Range("E46").Value = '  +8.10%  '
$ws.Range("D47").Value = '0.05861'
$ws.Range("E47").Value = '  +2.73%  '
$ws.Range("D48").Value = '0.3941'
$ws.Range("E48").Value = '  +6.12%  '
$ws.Range("D49").Value = '8.574'
$ws.Range("E49").Value = '  +3.58%  '
$ws.Range("D50").Value = '33.44'
$ws.Range("E50").Value = '  +3.04%  '
$ws.Range("E51").Value = '  +4.00%  '
